$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12, shifting rows 12:57 down to 13:58
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with data (copy pattern from the original row-12 data,
# now shifted to row 13, but with the updated fields from the diff)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 45030
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107001
$ws.Range("J12").Value = "Caqui"
$ws.Range("K12").Value = "Fuyu"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 23000
$ws.Range("P12").Value = 22455
$ws.Range("Q12").Value = "$/bandeja 15 kilos granel"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1497
$ws.Range("T12").Value = 15
